$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.108.02"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "2.059.96"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'230.53"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +6.37%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "2.366.00"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "'14.61"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "'20.67"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "'0.754"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "2.062.67"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "37.997.28"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'69.94"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").Value = "'224.85"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'2.46"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'166.50"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  +6.49%  "
$ws.Range("D29").Value = "'19.03"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "'1.36"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'1.99"
$ws.Range("E35").Value = "  +7.81%  "
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = "  +13.68%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  +4.94%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'98.50"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0219"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "1.480.84"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'0.0941"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.86"
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'16.69"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'4.11"
$ws.Range("E47").Value = "  +16.09%  "
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'2.96"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").Value = "2.255.22"
$ws.Range("E51").Value = "  +2.58%  "
